$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Colorado's sales territory moved from the Central Region to the Western Region.
# Before the re-sort below, Colorado is on row 7 (State/StateFullName/SalesRegion).
$ws.Range("C7").Value = "Western Region"

# Re-sort the SalesRegions table alphabetically by State (column A), ascending,
# with the header row excluded from the sort.
$lo = $ws.ListObjects.Item("SalesRegions")
$lo.Sort.SortFields.Clear()
$lo.Sort.SortFields.Add2($ws.Range("A2:A51"), 0, 1)
$lo.Sort.Header = 1
$lo.Sort.Apply()

# Leave the selection on C5, where Excel landed after the sort/edit.
$ws.Range("C5").Select() | Out-Null
